$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.315, 0.169, 0.323, 16, 14, -60),
    @(-0.322, -0.014, 0.032, 15, 13, 112),
    @(-0.297, 0.07099999999999999, 0.045, -68, 26, 104),
    @(0.012, 0.221, 0.043, -87, -16, -106),
    @(0.178, 0.118, 0.055, -3, -77, -134),
    @(0.236, 0.179, 0.355, -126, -15, 14),
    @(-0.07099999999999999, -0.095, 0.044, 38, -43, -131),
    @(0.106, -0.209, 0.078, 68, -4, 68),
    @(0.014, -0.07000000000000001, 0.305, -176, -29, 77),
    @(-0.131, -0.237, 0.525, -73, 7, -139)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
